$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newIds = @{
    2 = "67e72c4acde4e894925703ff"
    3 = "67e72c4acde4e894925703fd"
    4 = "67e72c4acde4e89492570401"
    5 = "67e72c49cde4e8949257024b"
    6 = "67e72c49cde4e8949257024f"
    7 = "67e72c49cde4e8949257024d"
    8 = "67e72c49cde4e89492570251"
    9 = "67e72c4acde4e89492570515"
    10 = "67e72c4acde4e89492570517"
    11 = "67e72c4acde4e89492570519"
    12 = "67e72c4acde4e89492570667"
    13 = "67e72c4acde4e89492570665"
    14 = "67e72c49cde4e894925703d0"
    15 = "67e72c49cde4e894925703d4"
    16 = "67e72c49cde4e894925703ce"
    17 = "67e72c49cde4e894925703d2"
    18 = "67e72c4acde4e8949257056f"
    19 = "67e72c4acde4e89492570571"
    20 = "67e72c49cde4e894925702cd"
    21 = "67e72c49cde4e894925702cb"
    22 = "67e72c49cde4e894925702c9"
    23 = "67e72c48cde4e894925701fb"
    24 = "67e72c48cde4e894925701fd"
    25 = "67e72c48cde4e894925701ff"
    26 = "67e72c49cde4e894925703a1"
    27 = "67e72c49cde4e8949257039f"
    28 = "67e72c49cde4e894925703a5"
    29 = "67e72c4acde4e8949257045b"
    30 = "67e72c4acde4e8949257045f"
    31 = "67e72c4acde4e89492570465"
    32 = "67e72c4acde4e8949257045d"
    33 = "67e72c4acde4e89492570461"
    34 = "67e72c49cde4e8949257028b"
    35 = "67e72c49cde4e8949257028d"
    36 = "67e72c49cde4e8949257028d"
    37 = "67e72c4acde4e89492570428"
    38 = "67e72c4acde4e89492570424"
    39 = "67e72c4acde4e89492570420"
    40 = "67e72c4acde4e89492570422"
    41 = "67e72c4acde4e89492570426"
    42 = "67e72c4acde4e8949257053a"
    43 = "67e72c4acde4e89492570536"
    44 = "67e72c4acde4e89492570538"
    45 = "67e72c4acde4e89492570534"
    46 = "67e72c4acde4e894925704a3"
    47 = "67e72c4acde4e894925704a5"
    48 = "67e72c4acde4e894925704a9"
    49 = "67e72c4acde4e894925704a7"
    50 = "67e72c49cde4e89492570276"
    51 = "67e72c49cde4e89492570278"
    52 = "67e72c4acde4e8949257067a"
    53 = "67e72c4acde4e8949257067c"
    54 = "67e72c4acde4e8949257067e"
    55 = "67e72c4acde4e894925704e9"
    56 = "67e72c4acde4e894925704e5"
    57 = "67e72c4acde4e894925704e3"
    58 = "67e72c4acde4e894925705c5"
    59 = "67e72c4acde4e894925705c3"
    60 = "67e72c4acde4e894925705c7"
    61 = "67e72c4acde4e894925705c9"
    62 = "67e72c49cde4e89492570317"
    63 = "67e72c49cde4e89492570319"
    64 = "67e72c49cde4e89492570315"
    65 = "67e72c49cde4e89492570351"
    66 = "67e72c49cde4e89492570355"
    67 = "67e72c49cde4e8949257033a"
    68 = "67e72c49cde4e89492570338"
    69 = "67e72c4acde4e894925705fe"
    70 = "67e72c4acde4e894925705fa"
    71 = "67e72c4acde4e894925705f6"
    72 = "67e72c49cde4e89492570372"
    73 = "67e72c49cde4e8949257036e"
    74 = "67e72c49cde4e89492570374"
    75 = "67e72c49cde4e89492570370"
    76 = "67e72c4bcde4e894925706ee"
    77 = "67e72c4bcde4e894925706f2"
    78 = "67e72c4bcde4e894925706f0"
    79 = "67e72c4bcde4e894925706f4"
    80 = "67e72c49cde4e894925702fd"
    81 = "67e72c4bcde4e8949257071b"
    82 = "67e72c4bcde4e89492570719"
    83 = "67e72c4bcde4e8949257071d"
    84 = "67e72c4acde4e89492570588"
    85 = "67e72c4acde4e8949257058a"
    86 = "67e72c4acde4e894925705a1"
    87 = "67e72c4acde4e89492570699"
    88 = "67e72c4acde4e894925705b8"
    89 = "67e72c49cde4e894925702b0"
    90 = "67e72c49cde4e894925702ae"
    91 = "67e72c49cde4e89492570222"
    92 = "67e72c49cde4e89492570224"
    93 = "67e72c49cde4e8949257021e"
    94 = "67e72c49cde4e89492570220"
    95 = "67e72c4acde4e89492570636"
    96 = "67e72c4acde4e89492570638"
    97 = "67e72c4acde4e8949257063c"
    98 = "67e72c4acde4e8949257063a"
    99 = "67e72c4acde4e89492570623"
    100 = "67e72c4acde4e89492570621"
    101 = "67e72c4acde4e89492570504"
    102 = "67e72c4bcde4e894925706b9"
    103 = "67e72c4bcde4e894925706c1"
    104 = "67e72c4bcde4e894925706c1"
    105 = "67e72c4bcde4e894925706bb"
    106 = "67e72c4bcde4e894925706bd"
    107 = "67e72c4bcde4e894925706a8"
    108 = "67e72c4bcde4e89492570746"
    109 = "67e72c4bcde4e89492570744"
    110 = "67e72c4bcde4e8949257078d"
    111 = "67e72c4bcde4e89492570791"
    113 = "67e72c4bcde4e8949257086f"
    115 = "67e72c4bcde4e894925707a6"
    116 = "67e72c4bcde4e89492570806"
    121 = "67e72c4bcde4e89492570823"
    122 = "67e72c4bcde4e89492570811"
    127 = "67e72c4bcde4e89492570833"
    128 = "67e72c4bcde4e894925707c8"
    129 = "67e72c4bcde4e89492570845"
    130 = "67e72c4bcde4e89492570852"
    131 = "67e72c4bcde4e89492570778"
    132 = "67e72c4bcde4e89492570862"
    133 = "67e72c4bcde4e894925707ee"
}

foreach ($row in $newIds.Keys) {
    $ws.Cells.Item($row, 22).Value = $newIds[$row]
}
